# agregar picker al mapa
# Update the "ubicacion" column (G) from "Concepcion" to "Santiago, Chile"
# for every data row, then move the active selection to B2 (picker reset).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Concepcion") {
        $cell.Value = "Santiago, Chile"
    }
}

# Re-fit the "ubicacion" / "modelo" columns now that the longer
# "Santiago, Chile" label lives in column G (and column H gets its own
# explicit width now too).
$ws.Columns.Item(7).ColumnWidth = 13
$ws.Columns.Item(8).ColumnWidth = 9.7

$ws.Range("B2").Select()
